# issue #5: property land done
#
# Normalises whitespace/punctuation in several existing text cells across
# the "土地" (land), "建物" (building) and "債務" (debt) sheets, and turns the
# "土地" sheet into a tidy-data table by renaming its header row to the
# generic English column names and appending metadata columns
# (property_category, category, date, legislator_name, legislator_id,
# source_file, index) to every data row.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet 1: 土地 (land)
# ---------------------------------------------------------------------
$land = $wb.Worksheets.Item(1)

# Header row -> generic tidy-data column names.
$land.Range("B1").Value = "name"
$land.Range("C1").Value = "area"
$land.Range("D1").Value = "share_portion"
$land.Range("E1").Value = "owner"
$land.Range("F1").Value = "register_date"
$land.Range("G1").Value = "register_reason"
$land.Range("H1").Value = "acquire_value"

# New trailing metadata headers, styled to match the existing bold/bordered
# header cells (B1:H1).
$newHeaders = @("I1","J1","K1","L1","M1","N1","O1")
$newHeaderNames = @("property_category","category","date","legislator_name","legislator_id","source_file","index")
for ($i = 0; $i -lt $newHeaders.Length; $i++) {
    $cell = $land.Range($newHeaders[$i])
    $cell.Value = $newHeaderNames[$i]
    $cell.Font.Bold = $true
    $cell.Borders.LineStyle = 1
    $cell.HorizontalAlignment = -4108
    $cell.VerticalAlignment = -4160
}

# Clean up whitespace / dash punctuation in existing data cells.
$land.Range("B2").Value = "臺中市大里區大孝段00380001地號"
$land.Range("F2").Value = "93年08月26日"
$land.Range("B3").Value = "臺中市大里區大孝段00380000地號"
$land.Range("F3").Value = "93年08月26日"

# Append the new metadata columns to both data rows. "date" (column K) is
# forced to a text number-format first so the ISO-looking "2013-12-30"
# string is kept literal instead of being auto-converted to a date serial.
$land.Range("K2:K3").NumberFormat = "@"

$land.Range("I2").Value = "land"
$land.Range("J2").Value = "normal"
$land.Range("K2").Value = "2013-12-30"
$land.Range("L2").Value = "何欣純"
$land.Range("M2").Value = 1733
$land.Range("N2").Value = "tmp8e3c1"
$land.Range("O2").Value = 14

$land.Range("I3").Value = "land"
$land.Range("J3").Value = "normal"
$land.Range("K3").Value = "2013-12-30"
$land.Range("L3").Value = "何欣純"
$land.Range("M3").Value = 1733
$land.Range("N3").Value = "tmp8e3c1"
$land.Range("O3").Value = 15

# ---------------------------------------------------------------------
# Sheet 2: 建物 (building) -- whitespace / dash clean-up only
# ---------------------------------------------------------------------
$building = $wb.Worksheets.Item(2)
$building.Range("B2").Value = "臺中市清水區秀水段秀水小段00060000建號"
$building.Range("F2").Value = "89年01月14日"
$building.Range("B3").Value = "臺中市清水區秀水段秀水小段01498000建號"
$building.Range("F3").Value = "93年08月26日"

# ---------------------------------------------------------------------
# Sheet 4: 債務 (debt) -- whitespace / punctuation clean-up only
# ---------------------------------------------------------------------
$debt = $wb.Worksheets.Item(4)
# Balance column ("餘額") keeps its original text type (was already a
# shared-string cell, e.g. "5,215,377") -- force text format so the now
# purely-numeric-looking "5215377"/"1500000" strings aren't reinterpreted
# as numbers.
$debt.Range("E2:E3").NumberFormat = "@"

$debt.Range("D2").Value = "霧峰鄉農會臺中市霧峰區四德路10號"
$debt.Range("E2").Value = "5215377"
$debt.Range("F2").Value = "93年09月14日"
$debt.Range("D3").Value = "台中商業銀行清水分行臺中市清水區中山路104號"
$debt.Range("E3").Value = "1500000"
$debt.Range("F3").Value = "102年09月26日"
